# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Acelga"
# as row 12, pushing the existing rows 12..132 down to 13..133.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 12 (shifts rows 12-132 -> 13-133)
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with the new record's data
$ws.Cells.Item(12, 1).Value = 4
$ws.Cells.Item(12, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(12, 3).Value = "Los Lagos"
$ws.Cells.Item(12, 4).Value = 44552
$ws.Cells.Item(12, 5).Value = 10
$ws.Cells.Item(12, 6).Value = 100112009
$ws.Cells.Item(12, 7).Value = "Acelga"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 20
$ws.Cells.Item(12, 11).Value = 10000
$ws.Cells.Item(12, 12).Value = 10000
$ws.Cells.Item(12, 13).Value = 10000
$ws.Cells.Item(12, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(12, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(12, 16).Value = 833
$ws.Cells.Item(12, 17).Value = 12
$ws.Cells.Item(12, 18).Value = "Hortaliza"
